$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated metric values (columns C=R^2, D=RMSE, E=U) for rows 2-10 ---
$ws.Range("C2").Value = -12.9743
$ws.Range("D2").Value = 1.3727
$ws.Range("E2").Value = 3.4134

$ws.Range("C3").Value = -6.204
$ws.Range("D3").Value = 1.2395
$ws.Range("E3").Value = 2.6942

$ws.Range("C4").Value = -3.2452
$ws.Range("D4").Value = 1.1736
$ws.Range("E4").Value = 2.537

$ws.Range("C5").Value = -1.1614
$ws.Range("D5").Value = 1.0078
$ws.Range("E5").Value = 2.1834

$ws.Range("C6").Value = -1.0167
$ws.Range("D6").Value = 1.0774
$ws.Range("E6").Value = 2.3707

$ws.Range("C7").Value = -1.2619
$ws.Range("D7").Value = 1.1244
$ws.Range("E7").Value = 2.8553

$ws.Range("C8").Value = -1.3815
$ws.Range("D8").Value = 1.1626
$ws.Range("E8").Value = 2.9188

$ws.Range("C9").Value = -1.3579
$ws.Range("D9").Value = 1.1655
$ws.Range("E9").Value = 2.8769

$ws.Range("C10").Value = -1.4435
$ws.Range("D10").Value = 1.1953
$ws.Range("E10").Value = 2.9102

# --- Refreshed color-scale fills on the RMSE/U columns (D:E) ---
$ws.Range("D2").Interior.Color = 16121079   # F7FCF5
$ws.Range("E2").Interior.Color = 16121079   # F7FCF5

$ws.Range("D3").Interior.Color = 10410660   # A4DA9E
$ws.Range("E3").Interior.Color = 6665042    # 52B365

$ws.Range("D4").Interior.Color = 5018668    # 2C944C
$ws.Range("E4").Interior.Color = 5018668    # 2C944C

$ws.Range("D5").Interior.Color = 1786880    # 00441B
$ws.Range("E5").Interior.Color = 1786880    # 00441B

$ws.Range("D6").Interior.Color = 3767314    # 127C39
$ws.Range("E6").Interior.Color = 3240711    # 077331

$ws.Range("D7").Interior.Color = 5413939    # 339C52
$ws.Range("E7").Interior.Color = 8637572    # 84CC83

$ws.Range("D8").Interior.Color = 6796629    # 55B567
$ws.Range("E8").Interior.Color = 9754008    # 98D594

$ws.Range("D9").Interior.Color = 6862424    # 58B668
$ws.Range("E9").Interior.Color = 9031563    # 8BCF89

$ws.Range("D10").Interior.Color = 8046201   # 79C67A
$ws.Range("E10").Interior.Color = 9556885   # 95D391

# --- Font color needs to flip to the light swatch (F1F1F1) on the two darkest fills ---
$ws.Range("D5:E6").Font.Color = 15856113
